$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Lokakuu")
$ws.Activate()

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("C4").Value = 5

$ws.Range("C5").Select()

$wsDec = $wb.Worksheets.Item("Joulukuu")
$wsDec.Activate()
$excel.ActiveWindow.ScrollRow = 7
$wsDec.Range("J2").Select()

$ws.Activate()

